$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "requirements"

# Header row
$ws.Range("A1").Value = "UID"
$ws.Range("B1").Value = "STATEMENT"
$ws.Range("C1").Value = "PARENT"

# Data rows
$ws.Range("A2").Value = "REQ-001"
$ws.Range("B2").Value = "Statement #1"
$ws.Range("C2").Value = ""

$ws.Range("A3").Value = "REQ-002"
$ws.Range("B3").Value = "Statement #2"

$ws.Range("A4").Value = "REQ-003"
$ws.Range("B4").Value = "Statement #3"

# Hyperlinked parent references
$ws.Hyperlinks.Add($ws.Range("C3"), "", "'requirements'!A2", "", "REQ-001")
$ws.Hyperlinks.Add($ws.Range("C4"), "", "'requirements'!A3", "", "REQ-002")

# Turn the range into a table
$lo = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $ws.Range("A1:C4"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$lo.Name = "Table1"
$lo.TableStyle = "TableStyleMedium9"

# Column widths (approximate target look)
$ws.Columns.Item(1).ColumnWidth = 7.66
$ws.Columns.Item(2).ColumnWidth = 12.66
$ws.Columns.Item(3).ColumnWidth = 9.66

# Selection as left by the author
$ws.Range("E10").Select() | Out-Null
